$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.465488910675049
$ws.Range("B1").Value = 4.309579849243164
$ws.Range("C1").Value = 3.226606845855713
$ws.Range("D1").Value = 0.9240348339080811
$ws.Range("E1").Value = 0.4779215753078461
